$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update raw metric values (Precision/Recall columns B/C/D) ---
$ws.Range("C3").Value = 0.77869999999999995

$ws.Range("B4").Value = 0.81494659185409501
$ws.Range("C4").Value = 0.84060000000000001
$ws.Range("D4").Value = 0.81489999999999996

$ws.Range("B5").Value = 0.78291815519332797
$ws.Range("C5").Value = 0.78939999999999999
$ws.Range("D5").Value = 0.78290000000000004

$ws.Range("B6").Value = 0.79003560543060303
$ws.Range("C6").Value = 0.78849999999999998
$ws.Range("D6").Value = 0.79

$ws.Range("B7").Value = 0.72953736782073897
$ws.Range("C7").Value = 0.73870000000000002
$ws.Range("D7").Value = 0.72950000000000004

# --- Apply the numeric (0.0000) format to cells that now need it ---
$ws.Range("C3").NumberFormat = "0.0000"
$ws.Range("B4:C4").NumberFormat = "0.0000"
$ws.Range("B5:D5").NumberFormat = "0.0000"
$ws.Range("B6:D6").NumberFormat = "0.0000"

# --- F1 Score column: compute as the harmonic mean of Recall and Precision ---
$ws.Range("E2").Formula = "=2*(C2*D2)/(C2+D2)"
$ws.Range("E3:E7").Formula = "=2*(C3*D3)/(C3+D3)"
$ws.Range("E3:E7").NumberFormat = "0.0000"

# --- Restore cursor/selection position ---
$ws.Range("H13").Select() | Out-Null
